# certificado de retencion segun res AGIP 382
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fill in status for existing backlog rows that were left blank
$ws.Range("B61").Value = "no comenzado"
$ws.Range("B62").Value = "no comenzado"
$ws.Range("B63").Value = "no comenzado"

# Add the new backlog item
$ws.Range("A64").Value = "reportes de ot no estan funcionando"
$ws.Range("B64").Value = "no comenzado"

$ws.Range("B68").Select()
